$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rows ---
# Row 2: street_num 12 -> 41, clear house_num (E2)
$ws.Range("C2").Value = 41
$ws.Range("E2").ClearContents()

# Row 3: street_num 24 -> 25, house_num (E3) = "1A"
$ws.Range("C3").Value = 25
$ws.Range("E3").Value = "1A"

# --- Add new row 5 (new sample record) ---
$ws.Range("A5").Value = "Mary Antonette"
$ws.Range("B5").Value = "ASDJ5612GJ"
$ws.Range("C5").Value = 47
$ws.Range("D5").Value = "80th Street East"
$ws.Range("F5").Value = "City centre"
$ws.Range("G5").Value = "Vladivostock"
$ws.Range("H5").Value = "Russia"
$ws.Range("I5").Value = 124654
$ws.Range("J5").Value = "gfkuavlk2"

# --- Update selection to match the saved state ---
$ws.Range("J5").Select()
